# dasde - excel formula conversion.
# Fill in rows 9, 10 and 11 with new redial tracking data (same pattern as
# the existing rows above), continuing the CONCATENATE/COUNTIF formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---------------------------------------------------------------
$ws.Range("A9").Value = "y"

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 44593

$ws.Range("C9").Value = 1234567

$ws.Range("D9").Formula = "=CONCATENATE(B9,C9)"
$ws.Range("E9").Formula = '=+IF(A9="n","",COUNTIF(D10:$D10007,$D9))'

$ws.Range("F2").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Formula = "=B9+1"

$ws.Range("G9").Formula = "=CONCATENATE(F9,C9)"
$ws.Range("H9").Formula = '=+IF(A9="n","",COUNTIF(D10:$D10007,$G9))'

$ws.Range("I2").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Formula = "=B9+2"

$ws.Range("J9").Formula = "=CONCATENATE(I9,C9)"
$ws.Range("K9").Formula = '=+IF(A9="n","",COUNTIF(D10:$D10007,$J9))'

$ws.Range("L2").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Row 10 ----------------------------------------------------------------
$ws.Range("A10").Value = "y"

$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 44593

$ws.Range("C10").Value = 1234567

$ws.Range("D10").Formula = "=CONCATENATE(B10,C10)"
$ws.Range("E10").Formula = '=+IF(A10="n","",COUNTIF(D11:$D10008,$D10))'

$ws.Range("F2").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Formula = "=B10+1"

$ws.Range("G10").Formula = "=CONCATENATE(F10,C10)"
$ws.Range("H10").Formula = '=+IF(A10="n","",COUNTIF(D11:$D10008,$G10))'

$ws.Range("I2").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Formula = "=B10+2"

$ws.Range("J10").Formula = "=CONCATENATE(I10,C10)"
$ws.Range("K10").Formula = '=+IF(A10="n","",COUNTIF(D11:$D10008,$J10))'

$ws.Range("L2").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- Row 11 ----------------------------------------------------------------
$ws.Range("A11").Value = "y"

$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = 44594

$ws.Range("C11").Value = 1234567

$ws.Range("D11").Formula = "=CONCATENATE(B11,C11)"
$ws.Range("E11").Formula = '=+IF(A11="n","",COUNTIF(D12:$D10009,$D11))'

$ws.Range("F2").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Formula = "=B11+1"

$ws.Range("G11").Formula = "=CONCATENATE(F11,C12)"
$ws.Range("H11").Formula = '=+IF(A11="n","",COUNTIF(D12:$D10009,$G11))'

$ws.Range("I2").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Formula = "=B11+2"

$ws.Range("J11").Formula = "=CONCATENATE(I11,C11)"
$ws.Range("K11").Formula = '=+IF(A11="n","",COUNTIF(D12:$D10009,$J11))'

# --- Selection (moved by the user to the new row) --------------------------
$ws.Range("A11:XFD11").Select() | Out-Null
